$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E3").Value = 16.5154
$ws.Range("C12").Value = -11.5925
$ws.Range("E14").Value = 16.608
$ws.Range("E26").Value = 16.52179999999999
$ws.Range("C27").Value = -13.2457
$ws.Range("E31").Value = 16.7146
$ws.Range("C32").Value = -13.6961
$ws.Range("E35").Value = 16.7298
$ws.Range("C36").Value = -13.35710000000002
$ws.Range("E37").Value = 16.77400000000001
$ws.Range("C38").Value = -13.94729999999999
$ws.Range("E45").Value = 16.5344
$ws.Range("C46").Value = -14.49799999999999
$ws.Range("E52").Value = 17.0662
$ws.Range("C54").Value = -13.5743
$ws.Range("C55").Value = -13.4034
$ws.Range("C56").Value = -12.35249999999999
$ws.Range("E57").Value = 16.8087
$ws.Range("C67").Value = -10.67650000000001
$ws.Range("C69").Value = -11.8492
$ws.Range("C72").Value = -11.6981
$ws.Range("E81").Value = 16.49209999999999
$ws.Range("C83").Value = -14.0351
$ws.Range("E83").Value = 16.3771
$ws.Range("C86").Value = -13.6632
$ws.Range("C91").Value = -10.1111
$ws.Range("C93").Value = -10.77070000000001
$ws.Range("C99").Value = -13.9787
$ws.Range("E100").Value = 16.3603
$ws.Range("E102").Value = 16.8973
